$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 10783.667
$ws.Range("J9").Value = 1140.4
$ws.Range("L9").Value = 1140.4
$ws.Range("N9").Value = -1478.4
$ws.Range("H38").Value = 1284
$ws.Range("I38").Value = 1284
$ws.Range("K38").Value = 3852
$ws.Range("M38").Value = -3480
$ws.Range("H46").Value = 3000
$ws.Range("I46").Value = 3000
$ws.Range("K46").Value = 9000
$ws.Range("M46").Value = -8881
$ws.Range("H60").Value = 3000
$ws.Range("I60").Value = 3000
$ws.Range("K60").Value = 9000
$ws.Range("M60").Value = -8516
$ws.Range("H76").Value = 4451.4707
$ws.Range("I76").Value = 3619.7
$ws.Range("K76").Value = 3619.7
$ws.Range("M76").Value = -3304.7
$ws.Range("H79").Value = 4451.4707
$ws.Range("I79").Value = 3619.7
$ws.Range("K79").Value = 3619.7
$ws.Range("M79").Value = -2527.7
$ws.Range("H113").Value = 6986.6875
$ws.Range("I113").Value = 4328.7
$ws.Range("J113").Value = 11416.667
$ws.Range("K113").Value = 4328.7
$ws.Range("L113").Value = 11416.667
$ws.Range("M113").Value = -1074.7
$ws.Range("N113").Value = -17924.667
$ws.Range("H136").Value = 71875
$ws.Range("J136").Value = 71875
$ws.Range("L136").Value = 71875
$ws.Range("N136").Value = -82075
$ws.Range("H138").Value = 10207614
$ws.Range("I138").Value = 1258
$ws.Range("J138").Value = 14710419
$ws.Range("K138").Value = 3774
$ws.Range("L138").Value = 44131257
$ws.Range("M138").Value = 1366
$ws.Range("N138").Value = -44141537

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 12127
$ws.Range("I19").Value = 12127
$ws.Range("K19").Value = 12127
$ws.Range("M19").Value = -11898
$ws.Range("H45").Value = 9541.625
$ws.Range("I45").Value = 9377.733
$ws.Range("K45").Value = 9377.733
$ws.Range("M45").Value = -9000.733

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 10220.667
$ws.Range("I107").Value = 2213.6
$ws.Range("K107").Value = 2213.6
$ws.Range("M107").Value = -293.5999999999999
$ws.Range("H135").Value = 56944.5
$ws.Range("J135").Value = 56944.5
$ws.Range("L135").Value = 56944.5
$ws.Range("N135").Value = -67084.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1721.75
$ws.Range("I16").Value = 1849.5
$ws.Range("J16").Value = 1594
$ws.Range("K16").Value = 1849.5
$ws.Range("L16").Value = 1594
$ws.Range("M16").Value = -1562.5
$ws.Range("N16").Value = -2168
$ws.Range("H58").Value = 2107.9092
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H99").Value = 4729.091
$ws.Range("I99").Value = 4540.125
$ws.Range("J99").Value = 5233
$ws.Range("K99").Value = 4540.125
$ws.Range("L99").Value = 5233
$ws.Range("M99").Value = -3042.125
$ws.Range("N99").Value = -8229
$ws.Range("H103").Value = 6062.273
$ws.Range("I103").Value = 6062.273
$ws.Range("K103").Value = 6062.273
$ws.Range("M103").Value = -4890.273
$ws.Range("H105").Value = 1956.8077
$ws.Range("I105").Value = 1037.2307
$ws.Range("J105").Value = 2876.3845
$ws.Range("K105").Value = 1037.2307
$ws.Range("L105").Value = 2876.3845
$ws.Range("M105").Value = 709.7692999999999
$ws.Range("N105").Value = -6370.3845
$ws.Range("H113").Value = 1721.75
$ws.Range("I113").Value = 1849.5
$ws.Range("J113").Value = 1594
$ws.Range("K113").Value = 1849.5
$ws.Range("L113").Value = 1594
$ws.Range("M113").Value = 320.5
$ws.Range("N113").Value = -5934
$ws.Range("H126").Value = 4729.091
$ws.Range("I126").Value = 4540.125
$ws.Range("J126").Value = 5233
$ws.Range("K126").Value = 13620.375
$ws.Range("L126").Value = 15699
$ws.Range("M126").Value = -11150.375
$ws.Range("N126").Value = -20639
$ws.Range("H134").Value = 25846.111
$ws.Range("I134").Value = 7829.3687
$ws.Range("J134").Value = 68635.875
$ws.Range("K134").Value = 23488.1061
$ws.Range("L134").Value = 205907.625
$ws.Range("M134").Value = -20953.1061
$ws.Range("N134").Value = -210977.625
$ws.Range("H136").Value = 2107.9092
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 854.38464
$ws.Range("I5").Value = 778.1667
$ws.Range("K5").Value = 2334.5001
$ws.Range("M5").Value = -2222.5001
$ws.Range("H14").Value = 87846.22
$ws.Range("I14").Value = 87846.22
$ws.Range("K14").Value = 263538.66
$ws.Range("M14").Value = -263365.66
$ws.Range("H23").Value = 292.5625
$ws.Range("J23").Value = 348.27274
$ws.Range("L23").Value = 1044.81822
$ws.Range("N23").Value = -1514.81822
$ws.Range("H97").Value = 1757.1666
$ws.Range("I97").Value = 1997.5
$ws.Range("K97").Value = 5992.5
$ws.Range("M97").Value = -5496.5
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H135").Value = 854.38464
$ws.Range("I135").Value = 778.1667
$ws.Range("K135").Value = 7003.5003
$ws.Range("M135").Value = -4468.5003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2590
$ws.Range("I80").Value = 2498.3333
$ws.Range("K80").Value = 2498.3333
$ws.Range("M80").Value = -1500.3333
$ws.Range("H83").Value = 2590
$ws.Range("I83").Value = 2498.3333
$ws.Range("K83").Value = 12491.6665
$ws.Range("M83").Value = -7499.666499999999
$ws.Range("H102").Value = 32263008
$ws.Range("I102").Value = 1860.4445
$ws.Range("K102").Value = 1860.4445
$ws.Range("M102").Value = -238.4445000000001
$ws.Range("H122").Value = 2912.0557
$ws.Range("J122").Value = 2681.6667
$ws.Range("L122").Value = 8045.000100000001
$ws.Range("N122").Value = -12945.0001
$ws.Range("H126").Value = 19651.928
$ws.Range("I126").Value = 20779
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 62337
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -59867
$ws.Range("N126").Value = -19940
$ws.Range("H141").Value = 111603.5
$ws.Range("J141").Value = 111603.5
$ws.Range("L141").Value = 111603.5
$ws.Range("N141").Value = -121963.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6640.846
$ws.Range("I7").Value = 6393.4287
$ws.Range("J7").Value = 6929.5
$ws.Range("K7").Value = 6393.4287
$ws.Range("L7").Value = 6929.5
$ws.Range("M7").Value = -6281.4287
$ws.Range("N7").Value = -7153.5
$ws.Range("H40").Value = 6891.9165
$ws.Range("I40").Value = 6001
$ws.Range("J40").Value = 7337.375
$ws.Range("K40").Value = 6001
$ws.Range("L40").Value = 7337.375
$ws.Range("M40").Value = -5865
$ws.Range("N40").Value = -7609.375
$ws.Range("H42").Value = 6999.6
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H49").Value = 6999.6
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H100").Value = 6800.4
$ws.Range("I100").Value = 2063.125
$ws.Range("J100").Value = 25749.5
$ws.Range("K100").Value = 2063.125
$ws.Range("L100").Value = 25749.5
$ws.Range("M100").Value = -1522.125
$ws.Range("N100").Value = -26831.5
$ws.Range("H118").Value = 5999
$ws.Range("J118").Value = 5999
$ws.Range("L118").Value = 5999
$ws.Range("N118").Value = -9313
$ws.Range("H126").Value = 6640.846
$ws.Range("I126").Value = 6393.4287
$ws.Range("J126").Value = 6929.5
$ws.Range("K126").Value = 19180.2861
$ws.Range("L126").Value = 20788.5
$ws.Range("M126").Value = -16710.2861
$ws.Range("N126").Value = -25728.5
$ws.Range("H136").Value = 5193.9375
$ws.Range("I136").Value = 5090.3
$ws.Range("J136").Value = 5366.6665
$ws.Range("K136").Value = 15270.9
$ws.Range("L136").Value = 16099.9995
$ws.Range("M136").Value = -12720.9
$ws.Range("N136").Value = -21199.9995
